$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new data rows for mouse 1337 (mirrors the pattern used for mouse 1329
# in rows 14-16: only the first row carries the mouse id, each row carries
# the area label, no date/depth supplied).
$ws.Range("A18").Value = 1337
$ws.Range("C18").Value = "V1"

$ws.Range("C19").Value = "LM"

$ws.Range("C20").Value = "LI"

# Move the active selection, matching the post-edit cursor position.
$ws.Range("I16").Select()
